$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header cell in A1 held the shared string "Grade"; change it to lowercase "grade".
$ws.Range("A1").Value = "grade"
